$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-6) are being re-shuffled: the Fecha (D), Volumen (J),
# Precio minimo (K), Precio maximo (L), Precio promedio ponderado (M) and
# Precio $/Kg (P) columns get a new set of values per row, while all other
# columns (A,B,C,E,F,G,H,I,N,O,Q,R) stay the same for every row.

$rows = @{
    2 = @{ D = 44638; J = 800;  K = 2500; L = 2800; M = 2650; P = 2650 }
    3 = @{ D = 44210; J = 1450; K = 1600; L = 1700; M = 1650; P = 1650 }
    4 = @{ D = 44175; J = 1400; K = 1900; L = 2000; M = 1950; P = 1950 }
    5 = @{ D = 44200; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
    6 = @{ D = 44537; J = 800;  K = 1300; L = 1400; M = 1350; P = 1350 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
